$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 38/39 swap: VeChain and InternetComputer(DFINITY) switch places
$ws.Range("B38").Value = "VeChain"
$ws.Range("C38").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02353"
$ws.Range("E38").Value = "  +2.67%  "

$ws.Range("B39").Value = "InternetComputer(DFINITY)"
$ws.Range("C39").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.264"
$ws.Range("E39").Value = "  +2.21%  "

# Price (D) and Volume(1h) (E) updates for remaining rows
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.256.18"
$ws.Range("E2").Value = "  +1.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.803.62"
$ws.Range("E3").Value = "  +2.82%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "337.38"
$ws.Range("E5").Value = "  +0.44%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4667"
$ws.Range("E7").Value = "  +22.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3786"
$ws.Range("E8").Value = "  +10.62%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.10"
$ws.Range("E9").Value = "  -2.18%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07635"
$ws.Range("E10").Value = "  +5.63%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.151"
$ws.Range("E11").Value = "  +2.68%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.39"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.336"
$ws.Range("E14").Value = "  +2.78%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.459"
$ws.Range("E15").Value = "  +4.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.805.57"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("E17").Value = "  +3.47%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06742"
$ws.Range("E18").Value = "  +2.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "81.85"
$ws.Range("E19").Value = "  +3.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.42"
$ws.Range("E21").Value = "  +4.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.414"
$ws.Range("E22").Value = "  +3.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.224.66"
$ws.Range("E23").Value = "  +0.89%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.86"
$ws.Range("E24").Value = "  +1.56%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.412"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "20.76"
$ws.Range("E26").Value = "  +4.61%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.91"
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  +3.00%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.010.74"
$ws.Range("E29").Value = "  +3.00%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.85"
$ws.Range("E30").Value = "  +1.89%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.256"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.034"
$ws.Range("E32").Value = "  +0.25%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09609"
$ws.Range("E33").Value = "  +8.97%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.859"
$ws.Range("E34").Value = "  +0.34%  "
$ws.Range("E35").Value = "  +6.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06371"
$ws.Range("E36").Value = "  +3.24%  "
$ws.Range("E37").Value = "  -0.69%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6641"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.235"
$ws.Range("E41").Value = "  +1.80%  "
$ws.Range("E42").Value = "  -0.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.257"
$ws.Range("E43").Value = "  +3.96%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.18"
$ws.Range("E44").Value = "  +2.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9999"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6112"
$ws.Range("E46").Value = "  +0.84%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.836"
$ws.Range("E47").Value = "  -0.04%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "130.10"
$ws.Range("E48").Value = "  +2.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.035"
$ws.Range("E49").Value = "  +1.28%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07165"
$ws.Range("E50").Value = "  +2.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.177"
$ws.Range("E51").Value = "  +0.68%  "
